# Update the practice-table answers to match the new key (commit c986bee).
# Each table cell is addressed and rewritten positionally (Table.Cell(row, col))
# rather than via Find/Replace, because some new answers duplicate other
# original answers elsewhere in the table (e.g. "92÷9=10, 2"), and a global
# Find/Replace could clobber an already-updated cell.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "34÷7=4, 6"
$t.Cell(1, 2).Range.Text = "63÷4=15, 3"
$t.Cell(1, 3).Range.Text = "93÷7=13, 2"
$t.Cell(1, 4).Range.Text = "92÷9=10, 2"
$t.Cell(1, 5).Range.Text = "63÷6=10, 3"

# Row 5
$t.Cell(5, 1).Range.Text = "83÷4=20, 3"
$t.Cell(5, 2).Range.Text = "37÷2=18, 1"
$t.Cell(5, 3).Range.Text = "45÷5=9, 0"
$t.Cell(5, 4).Range.Text = "34÷2=17, 0"
$t.Cell(5, 5).Range.Text = "35÷8=4, 3"

# Row 9
$t.Cell(9, 1).Range.Text = "72÷9=8, 0"
$t.Cell(9, 2).Range.Text = "97÷3=32, 1"
$t.Cell(9, 3).Range.Text = "24÷7=3, 3"
$t.Cell(9, 4).Range.Text = "47÷7=6, 5"
$t.Cell(9, 5).Range.Text = "42÷7=6, 0"

# Row 13
$t.Cell(13, 1).Range.Text = "77÷9=8, 5"
$t.Cell(13, 2).Range.Text = "52÷3=17, 1"
$t.Cell(13, 3).Range.Text = "81÷6=13, 3"
$t.Cell(13, 4).Range.Text = "69÷8=8, 5"
$t.Cell(13, 5).Range.Text = "26÷3=8, 2"

# Row 17
$t.Cell(17, 1).Range.Text = "37÷6=6, 1"
$t.Cell(17, 2).Range.Text = "10÷4=2, 2"
$t.Cell(17, 3).Range.Text = "33÷5=6, 3"
$t.Cell(17, 4).Range.Text = "98÷7=14, 0"
$t.Cell(17, 5).Range.Text = "85÷6=14, 1"
